$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 192
$ws.Range("F4").Value = 1141
$ws.Range("F5").Value = 41
$ws.Range("F10").Value = 339
$ws.Range("F11").Value = 420
$ws.Range("F14").Value = 354
$ws.Range("F17").Value = 474
$ws.Range("F18").Value = 446
$ws.Range("F19").Value = 5601
$ws.Range("F21").Value = 1563
$ws.Range("F23").Value = 4780
$ws.Range("F26").Value = 1502
$ws.Range("F29").Value = 653
$ws.Range("F30").Value = 66
$ws.Range("F32").Value = 3794

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 136
$ws.Range("F8").Value = 104

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9384
$ws.Range("F3").Value = 581
$ws.Range("F4").Value = 2131

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9384
$ws.Range("F3").Value = 581
$ws.Range("F4").Value = 2132
$ws.Range("F5").Value = 192
$ws.Range("F7").Value = 1141
$ws.Range("F8").Value = 41
$ws.Range("F12").Value = 339
$ws.Range("F13").Value = 420
$ws.Range("F16").Value = 354
$ws.Range("F22").Value = 446
$ws.Range("F23").Value = 5601
$ws.Range("F25").Value = 1563
$ws.Range("F31").Value = 4780
$ws.Range("F34").Value = 1502
$ws.Range("F37").Value = 653
$ws.Range("F38").Value = 66
$ws.Range("F46").Value = 3794
